# Update DataEntry proposal-content cells with the revised wording from the author's edit pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

$ws.Cells.Item(2, 2).Value = 'Sunlight Sensor'

$ws.Cells.Item(3, 2).Value = 'Raphael Carlo Najera'

$ws.Cells.Item(4, 2).Value = 'https://github.com/RaphaelNajera/Sunlight_Sensor'

$ws.Cells.Item(5, 2).Value = 'use the sunlight sensor to monitor sunlight intensity, IR intensity and UV intensity. This will give data for UV-Light, visible light and infrared light.'

$ws.Cells.Item(6, 2).Value = 'the data of UV-light, visible light and infrared light when it get the data  from detect sunlight. It will also record the time when it received the data. It will also store message and output the message on a screen. '

$ws.Cells.Item(7, 2).Value = 'displaying the total data of uv index, visible light (in Lumens) and infrared light (in Lumens). It will also show record on the past day and also display helpful message when it good to go outside and reporting helpful information like it is required to put sunscreen. '

$ws.Cells.Item(8, 2).Value = 'Humber College Institute of Technology & Advanced Learning North Campus Prototype Lab, Weather network, Raspberry Pi.'

$ws.Cells.Item(9, 2).Value = 'Johnson Liang and Adrian Caprini.'

$ws.Cells.Item(10, 2).Value = 'The sunlight sensor will dectect UV-light, visible light and infrared light. With the data we can measure the total visible light (in Lumens), infared light (in Lumens) and UV (UV index). Without it, the people will not know how long  they can stay outside till they get sun burn which damage your skin and can also cause skin cancer.'

$ws.Cells.Item(11, 2).Value = 'The sensor I''m going to use is Sunlight sensor. The sunlight sensor will monitor sunlight intensity, IR (Infrared light) intensity and UV (Ultraviolet light) intensity. I can use this sensor to detect and gathear data of UV-light (in uv index), visible light (in lumens) and infrared light (in Lemens). For example, if the uv light is 2, it means the uv light is low and if the uv light is 6, it mean the uv light is high. For data in lumens if you get the number 50, it means "cloudy day outdoor". With this I can output a helpful message to the user.'

$ws.Cells.Item(12, 2).Value = 'Simple IOT Sunlight Sensing Raspberry Pi Project - SunIOT Part 1. (2016, October 18). Retrieved September 17, 2017, from http://www.switchdoc.com/2016/10/simple-iot-sunlight-sensing-raspberry-pi-project-part-1/'

$ws.Cells.Item(13, 2).Value = 'Mazzillo, M., Shukla, P., & Mallik, R. (2010, September 27). 4H-SiC Schottky Photodiode Based Demonstrator Board for UV-Index Monitoring. Retrieved September 17, 2017, from http://ieeexplore.ieee.org/document/5585671/'

$ws.Cells.Item(14, 2).Value = 'Raspberry Pie 3: Use to connect the sunlight sensor to the main component. It will store the code to run the hardware and save data.

Pi2Grover - Grove Connector Interface for the Raspberry Pi: Provides the connection between Raspberry Pi pins and external Grove module.

Grove Sunlight  / IR / UV I2C sensor: Detect UV-light, visible light and infrared light.

'

$ws.Cells.Item(15, 2).Value = 'With the information from the sunlight sensor we can get data like UV index. This will help people to know the information of the sunlight each day. When it detect the uv light it will output a helpful message to the user. For example, if the uv is low it  the message would be "When doing outdoor activity minimal sun protection is required". If the uv is high the message would be "Sun protection required for outdoor activity".'

# Column A got narrower after the text edits (author resized it in Excel).
$ws.Columns.Item(1).ColumnWidth = 18.8

# Row heights were re-auto-sized by Excel after the wrapped text changed length.
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 75
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 120
$ws.Rows.Item(14).RowHeight = 135

# The saved selection moved to B10 with the view scrolled down to row 11.
$ws.Range("B10").Select()

# The DataBase summary sheet pulls everything via formulas; its row 2 was
# also re-auto-sized once the referenced text got shorter.
$ws2 = $wb.Worksheets.Item("DataBase")
$ws2.Rows.Item(2).RowHeight = 75

